$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

Set-TextValue "B2" "0.17"
Set-TextValue "B3" "-0.01"
Set-TextValue "B4" "-0.09"
Set-TextValue "C2" "44.29***"
Set-TextValue "C3" "2.21***"
Set-TextValue "C4" "0.98"
Set-TextValue "D2" "-0.89"
Set-TextValue "D3" "0.46***"
Set-TextValue "D4" "0.82*"
